$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting C:AG to D:AH.
$ws.Columns("C").Insert()

# Populate the newly inserted column with the "runCase" command block.
$ws.Range("C1").Value = "runCase"
$ws.Range("C2").Value = "xlsx,C:\projs\auto-test\testcases\testcase2.xlsx"
$ws.Range("C3").Value = 1

# Match column B's width/formatting (custom width, not best-fit).
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth
